$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2180.0833
$ws.Range("I6").Value = 80.5
$ws.Range("K6").Value = 241.5
$ws.Range("M6").Value = -129.5
$ws.Range("H8").Value = 1026
$ws.Range("H41").Value = 14285995
$ws.Range("I41").Value = 387.7143
$ws.Range("J41").Value = 28571604
$ws.Range("K41").Value = 387.7143
$ws.Range("L41").Value = 28571604
$ws.Range("M41").Value = 52.28570000000002
$ws.Range("N41").Value = -28572484
$ws.Range("H61").Value = 100093.7
$ws.Range("I61").Value = 117.125
$ws.Range("K61").Value = 351.375
$ws.Range("M61").Value = -179.375
$ws.Range("H82").Value = 6110.5
$ws.Range("H85").Value = 6110.5
$ws.Range("H111").Value = 869.8889
$ws.Range("I111").Value = 841
$ws.Range("J111").Value = 1101
$ws.Range("K111").Value = 2523
$ws.Range("L111").Value = 3303
$ws.Range("M111").Value = 544
$ws.Range("N111").Value = -9437
$ws.Range("H116").Value = 2977.8333
$ws.Range("I116").Value = 2655.25
$ws.Range("J116").Value = 3623
$ws.Range("K116").Value = 2655.25
$ws.Range("L116").Value = 3623
$ws.Range("M116").Value = 786.75
$ws.Range("N116").Value = -10507

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6283.864
$ws.Range("I61").Value = 5912.25
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 5912.25
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -5700.25
$ws.Range("N61").Value = -10424
$ws.Range("H132").Value = 1876.8462
$ws.Range("I132").Value = 1475.4762
$ws.Range("J132").Value = 3562.6
$ws.Range("K132").Value = 4426.4286
$ws.Range("L132").Value = 10687.8
$ws.Range("M132").Value = -1896.4286
$ws.Range("N132").Value = -15747.8
$ws.Range("H136").Value = 6283.864
$ws.Range("I136").Value = 5912.25
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 17736.75
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -15186.75
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 51593.332
$ws.Range("J74").Value = 51593.332
$ws.Range("L74").Value = 51593.332
$ws.Range("N74").Value = -53465.332
$ws.Range("H77").Value = 51593.332
$ws.Range("J77").Value = 51593.332
$ws.Range("L77").Value = 154779.996
$ws.Range("N77").Value = -164139.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 20075004
$ws.Range("I12").Value = 80000000
$ws.Range("J12").Value = 100006
$ws.Range("K12").Value = 80000000
$ws.Range("L12").Value = 100006
$ws.Range("M12").Value = -79999830
$ws.Range("N12").Value = -100346
$ws.Range("H16").Value = 1858.4117
$ws.Range("I16").Value = 1130
$ws.Range("J16").Value = 2899
$ws.Range("K16").Value = 1130
$ws.Range("L16").Value = 2899
$ws.Range("M16").Value = -843
$ws.Range("N16").Value = -3473
$ws.Range("H113").Value = 1858.4117
$ws.Range("I113").Value = 1130
$ws.Range("J113").Value = 2899
$ws.Range("K113").Value = 1130
$ws.Range("L113").Value = 2899
$ws.Range("M113").Value = 1040
$ws.Range("N113").Value = -7239

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1055.8889
$ws.Range("I10").Value = 180.25
$ws.Range("K10").Value = 540.75
$ws.Range("M10").Value = -401.75
$ws.Range("H13").Value = 750
$ws.Range("I13").Value = 333.33334
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 1000.00002
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = -832.0000200000001
$ws.Range("N13").Value = -6336
$ws.Range("H15").Value = 2383.3333
$ws.Range("J15").Value = 3500
$ws.Range("L15").Value = 10500
$ws.Range("N15").Value = -10780
$ws.Range("H47").Value = 1961.4445
$ws.Range("I47").Value = 83
$ws.Range("K47").Value = 249
$ws.Range("M47").Value = 182
$ws.Range("H64").Value = 10580.833
$ws.Range("I64").Value = 1794
$ws.Range("J64").Value = 16857.143
$ws.Range("K64").Value = 5382
$ws.Range("L64").Value = 50571.429
$ws.Range("M64").Value = -5112
$ws.Range("N64").Value = -51111.429
$ws.Range("H67").Value = 10580.833
$ws.Range("I67").Value = 1794
$ws.Range("J67").Value = 16857.143
$ws.Range("K67").Value = 5382
$ws.Range("L67").Value = 50571.429
$ws.Range("M67").Value = -4446
$ws.Range("N67").Value = -52443.429
$ws.Range("H70").Value = 2788.9092
$ws.Range("J70").Value = 3999.8
$ws.Range("L70").Value = 11999.4
$ws.Range("N70").Value = -12629.4
$ws.Range("H73").Value = 2788.9092
$ws.Range("J73").Value = 3999.8
$ws.Range("L73").Value = 11999.4
$ws.Range("N73").Value = -14183.4
$ws.Range("H76").Value = 3933
$ws.Range("I76").Value = 3799
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 11397
$ws.Range("L76").Value = 12000
$ws.Range("M76").Value = -11014
$ws.Range("N76").Value = -12766
$ws.Range("H79").Value = 3933
$ws.Range("I79").Value = 3799
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 11397
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = -10071
$ws.Range("N79").Value = -14652
$ws.Range("H82").Value = 1149.6666
$ws.Range("I82").Value = 1149.6666
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3448.9998
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3042.9998
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1149.6666
$ws.Range("I85").Value = 1149.6666
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3448.9998
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2044.9998
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 21281.2
$ws.Range("I19").Value = 2133.3333
$ws.Range("K19").Value = 2133.3333
$ws.Range("M19").Value = -1845.3333
$ws.Range("H113").Value = 2031.2
$ws.Range("I113").Value = 1626.8572
$ws.Range("J113").Value = 2385
$ws.Range("K113").Value = 1626.8572
$ws.Range("L113").Value = 2385
$ws.Range("M113").Value = 543.1428000000001
$ws.Range("N113").Value = -6725
$ws.Range("H125").Value = 19495
$ws.Range("J125").Value = 19495
$ws.Range("L125").Value = 19495
$ws.Range("N125").Value = -24415
$ws.Range("H132").Value = 3649.5676
$ws.Range("I132").Value = 3968.6667
$ws.Range("J132").Value = 3347.2632
$ws.Range("K132").Value = 11906.0001
$ws.Range("L132").Value = 10041.7896
$ws.Range("M132").Value = -9376.000100000001
$ws.Range("N132").Value = -15101.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4688.4546
$ws.Range("I31").Value = 6246.625
$ws.Range("J31").Value = 533.3333
$ws.Range("K31").Value = 6246.625
$ws.Range("L31").Value = 533.3333
$ws.Range("M31").Value = -5998.625
$ws.Range("N31").Value = -1029.3333
$ws.Range("H35").Value = 15140.25
$ws.Range("I35").Value = 280.5
$ws.Range("K35").Value = 280.5
$ws.Range("M35").Value = 55.5
$ws.Range("H45").Value = 7799.8
$ws.Range("I45").Value = 5999.5
$ws.Range("K45").Value = 5999.5
$ws.Range("M45").Value = -5592.5
$ws.Range("H46").Value = 1203.9656
$ws.Range("I46").Value = 917.1739
$ws.Range("J46").Value = 2303.3333
$ws.Range("K46").Value = 917.1739
$ws.Range("L46").Value = 2303.3333
$ws.Range("M46").Value = -729.1739
$ws.Range("N46").Value = -2679.3333
$ws.Range("H61").Value = 4307
$ws.Range("I61").Value = 1460.6666
$ws.Range("K61").Value = 1460.6666
$ws.Range("M61").Value = -1258.6666
$ws.Range("H113").Value = 4307
$ws.Range("I113").Value = 1460.6666
$ws.Range("K113").Value = 1460.6666
$ws.Range("M113").Value = 709.3334
$ws.Range("H124").Value = 29000
$ws.Range("J124").Value = 29000
$ws.Range("L124").Value = 29000
$ws.Range("N124").Value = -38820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5040002
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999886
$ws.Range("H58").Value = 13555.556
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9692
$ws.Range("H107").Value = 1298.7142
$ws.Range("I107").Value = 573.5
$ws.Range("J107").Value = 5650
$ws.Range("K107").Value = 1720.5
$ws.Range("L107").Value = 16950
$ws.Range("M107").Value = 199.5
$ws.Range("N107").Value = -20790

